# Daily attendance processing - 2026-01-07 10:07:19
# Reorders the comma-separated "Recorded By" values in column G:
#   - If the list contains an entry that is exactly "System" (capital S),
#     that entry is moved to the end of the list (the rest keep their
#     relative order).
#   - Otherwise (e.g. plain two-party lists with no "System" entry), the
#     list order is reversed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Length -le 1) {
        continue
    }

    $systemIndex = -1
    for ($i = 0; $i -lt $parts.Length; $i++) {
        if ($parts[$i] -eq "System") {
            $systemIndex = $i
        }
    }

    $newParts = @()
    if ($systemIndex -ge 0) {
        for ($i = 0; $i -lt $parts.Length; $i++) {
            if ($i -ne $systemIndex) {
                $newParts += $parts[$i]
            }
        }
        $newParts += "System"
    }
    else {
        for ($i = $parts.Length - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newText = $newParts -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
